# Parameters.xlsx update:
#  - swap the search keyword "covid19" -> "USA" (Sheet1 Keyword rows)
#  - lower the "Number of articles" target from 20 -> 5
#  - add a new "question number" parameters block in columns G:H
#    (Question 1/Paragraph, Question 2/Source)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keyword block (A1:E3) -------------------------------------------------
$ws.Range("B2").Value = "USA"
$ws.Range("E2").Value = "USA"
$ws.Range("B3").Value = 5
$ws.Range("E3").Value = 5

# New question-number parameters block (G1:H2) --------------------------
$ws.Range("G1").Value = "Question 1"
$ws.Range("G2").Value = "Question 2"
$ws.Range("H1").Value = "Paragraph"
$ws.Range("H2").Value = "Source"

# Widen the new column to fit its contents, same as the other columns
$ws.Columns("G").ColumnWidth = 9.67

# Leave the selection where the author left it when they saved the file
$ws.Range("N5").Select() | Out-Null
